$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 13) mirroring the existing table layout:
# Student Id | Trial | Correct | Elapsed Time | Date
$ws.Cells.Item(13, 1).Value = "ikleiman@stonybrook.edu"
$ws.Cells.Item(13, 2).Value = 2

# "false" and the date string would otherwise be auto-coerced by Excel into
# a Boolean / date serial number. Force literal text entry (leading
# apostrophe, like typing it in the UI) then clear the resulting
# quote-prefix formatting so the cell ends up as a plain text value with
# the default style, matching the rest of the sheet.
$c = $ws.Cells.Item(13, 3)
$c.Value = "'false"
$c.ClearFormats()

$ws.Cells.Item(13, 4).Value = 8

$e = $ws.Cells.Item(13, 5)
$e.Value = "'2019-12-30"
$e.ClearFormats()
